# Update the "paises" COVID-19 dashboard sheet with refreshed data
# and the new "last updated" timestamp, matching the upstream data
# refresh (countries & provincias Spain).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- 1. Update "last updated" banner text -------------------------------
$ws.Range("A1").Value = "Datos actualizados a 29 de Junio de 2020 a las 15:43"

# --- 2. Refresh numeric stats for countries whose totals changed --------
# Columns: B=Casos totales, C=Nuevos casos, D=Casos activos,
#          E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 2637439
$ws.Range("C4").Value = 362
$ws.Range("D4").Value = 1093612
$ws.Range("E4").Value = 1415384
$ws.Range("G4").Value = 6
$ws.Range("H4").Value = 128443

# Row 7 - India
$ws.Range("B7").Value = 554386
$ws.Range("C7").Value = 5189
$ws.Range("D7").Value = 324912
$ws.Range("E7").Value = 212906
$ws.Range("G7").Value = 81
$ws.Range("H7").Value = 16568

# Row 18 - Arabia Saudita
$ws.Range("B18").Value = 186436
$ws.Range("C18").Value = 3943
$ws.Range("D18").Value = 127118
$ws.Range("E18").Value = 57719
$ws.Range("G18").Value = 48
$ws.Range("H18").Value = 1599

# Row 29 - Belgica
$ws.Range("B29").Value = 61361
$ws.Range("C29").Value = 66
$ws.Range("E29").Value = 34688

# Row 30 - Argentina
$ws.Range("D30").Value = 21138
$ws.Range("E30").Value = 37550
$ws.Range("G30").Value = 13
$ws.Range("H30").Value = 1245

# Row 39 - Portugal
$ws.Range("B39").Value = 41912
$ws.Range("C39").Value = 266
$ws.Range("D39").Value = 27205
$ws.Range("E39").Value = 13139
$ws.Range("G39").Value = 4
$ws.Range("H39").Value = 1568

# Row 62
$ws.Range("B62").Value = 14288
$ws.Range("C62").Value = 242
$ws.Range("D62").Value = 12581
$ws.Range("E62").Value = 1433
$ws.Range("G62").Value = 4
$ws.Range("H62").Value = 274

# Row 72
$ws.Range("B72").Value = 8862
$ws.Range("C72").Value = 7
$ws.Range("E72").Value = 475

# Row 74
$ws.Range("E74").Value = 2750
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = 23

# Row 156 - Surinam
$ws.Range("B156").Value = 492
$ws.Range("C156").Value = 2
$ws.Range("E156").Value = 281
$ws.Range("G156").Value = 1
$ws.Range("H156").Value = 12

# Row 161 - Birmania
$ws.Range("D161").Value = 221
$ws.Range("E161").Value = 72

# --- 3. Siria overtakes Angola in the ranking ----------------------------
# Siria's refreshed totals are now higher than Angola's, so row 163 now
# shows Siria's (updated) figures and row 164 shows Angola's (prior,
# unchanged) figures - the two countries swap positions.
$ws.Range("A163").Value = "Siria"
$ws.Range("B163").Value = 269
$ws.Range("C163").Value = 13
$ws.Range("D163").Value = 102
$ws.Range("E163").Value = 158
$ws.Range("F163").Value = 0
$ws.Range("G163").Value = 0
$ws.Range("H163").Value = 9

$ws.Range("A164").Value = "Angola"
$ws.Range("B164").Value = 267
$ws.Range("C164").Value = 0
$ws.Range("D164").Value = 93
$ws.Range("E164").Value = 163
$ws.Range("F164").Value = 0
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 11
